# Update specific numeric values in Sheet1 as per target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D7").Value  = -7.364299999999991
$ws.Range("B10").Value = 8.609000000000002
$ws.Range("B12").Value = 6.076299999999999
$ws.Range("C13").Value = -12.62719999999998
$ws.Range("B18").Value = 4.797400000000008
$ws.Range("D20").Value = -8.463400000000002
